$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at G: this shifts the original "Regression" header / data
# column (G) to H, carrying its bestFit width along, and leaves a blank column G.
$ws.Columns.Item(7).Insert()

# New column G width (category column): 15.1666... -> stored xlsx width 16
$ws.Columns.Item(7).ColumnWidth = 15.166666666666666

# --- Brand-new TC id strings first, in order TC_02..TC_14 ---
$ws.Range("B3").Value  = "TC_02"
$ws.Range("B4").Value  = "TC_03"
$ws.Range("B5").Value  = "TC_04"
$ws.Range("B6").Value  = "TC_05"
$ws.Range("B7").Value  = "TC_06"
$ws.Range("B8").Value  = "TC_07"
$ws.Range("B9").Value  = "TC_08"
$ws.Range("B10").Value = "TC_09"
$ws.Range("B11").Value = "TC_10"
$ws.Range("B12").Value = "TC_11"
$ws.Range("B13").Value = "TC_12"
$ws.Range("B14").Value = "TC_13"
$ws.Range("B15").Value = "TC_14"

# --- Add new G column references that reuse existing strings before the old
#     F-column cells referencing them are overwritten, so "Smoke" and
#     "Regression" keep their original relative slot in the shared string table ---
$ws.Range("G2").Value  = "Regression"
$ws.Range("G3").Value  = "Regression"
$ws.Range("G4").Value  = "Smoke"
$ws.Range("G5").Value  = "Regression"

# --- Overwrite the old F2:F5 cells (previously Title 1 / Title 2 / Smoke /
#     Title 4) with the new Test Case Title strings, in order ---
$ws.Range("F2").Value = "Login"
$ws.Range("F3").Value = "Login with incorrect passoword"
$ws.Range("F4").Value = "Homepage"
$ws.Range("F5").Value = "Logo"

# --- Brand new title strings for rows 6 and 7 ---
$ws.Range("F6").Value = "Header"
$ws.Range("F7").Value = "Contact Us page"

# --- Remaining G column cells (rows 6-15), reusing existing strings ---
$ws.Range("G6").Value  = "Regression"
$ws.Range("G7").Value  = "Regression"
$ws.Range("G8").Value  = "Smoke"
$ws.Range("G9").Value  = "Smoke"
$ws.Range("G10").Value = "Regression"
$ws.Range("G11").Value = "Smoke"
$ws.Range("G12").Value = "Regression"
$ws.Range("G13").Value = "Smoke"
$ws.Range("G14").Value = "Regression"
$ws.Range("G15").Value = "Regression"

# --- Remaining F column cells for rows 8-15 (reusing Smoke/Regression) ---
$ws.Range("F8").Value  = "Smoke"
$ws.Range("F9").Value  = "Smoke"
$ws.Range("F10").Value = "Regression"
$ws.Range("F11").Value = "Smoke"
$ws.Range("F12").Value = "Regression"
$ws.Range("F13").Value = "Smoke"
$ws.Range("F14").Value = "Regression"
$ws.Range("F15").Value = "Regression"

# --- A / C / D / E / H columns for new rows 6-15 ---
for ($r = 6; $r -le 15; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
    $ws.Cells.Item($r, 3).Value = "Prathamesh"
    $ws.Cells.Item($r, 4).Value = "Lad"
    $ws.Cells.Item($r, 5).Value = "Chrome"
    $ws.Cells.Item($r, 8).Value = "Yes"
}

# --- Selection ---
$ws.Range("I25").Select()
